$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing the existing rows 151-180
# down to 152-181 (and extending the sheet dimension to A1:R181).
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record.
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44522
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 100112043
$ws.Range("G151").Value = "Pepino ensalada"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 150
$ws.Range("K151").Value = 11000
$ws.Range("L151").Value = 11000
$ws.Range("M151").Value = 11000
$ws.Range("N151").Value = "$/caja 60 unidades"
$ws.Range("O151").Value = "Región de Arica y Parinacota"
$ws.Range("P151").Value = 183
$ws.Range("Q151").Value = 60
$ws.Range("R151").Value = "Hortaliza"
